# Scheduled runner update: refresh market-price / profit figures on the
# per-job "Ridill_Profits" leve tables (columns H-N: currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 159098540
$ws.Cells.Item(62, 9).Value = 62510500
$ws.Cells.Item(62, 10).Value = 416666660
$ws.Cells.Item(62, 11).Value = 62510500
$ws.Cells.Item(62, 12).Value = 416666660
$ws.Cells.Item(62, 13).Value = -62509876
$ws.Cells.Item(62, 14).Value = -416667908

$ws.Cells.Item(65, 8).Value = 159098540
$ws.Cells.Item(65, 9).Value = 62510500
$ws.Cells.Item(65, 10).Value = 416666660
$ws.Cells.Item(65, 11).Value = 312552500
$ws.Cells.Item(65, 12).Value = 2083333300
$ws.Cells.Item(65, 13).Value = -312549380
$ws.Cells.Item(65, 14).Value = -2083339540

$ws.Cells.Item(98, 8).Value = 54731588
$ws.Cells.Item(98, 9).Value = 25000926
$ws.Cells.Item(98, 11).Value = 25000926
$ws.Cells.Item(98, 13).Value = -24999428

$ws.Cells.Item(111, 8).Value = 10307.077
$ws.Cells.Item(111, 9).Value = 3577.7778
$ws.Cells.Item(111, 10).Value = 25448
$ws.Cells.Item(111, 11).Value = 10733.3334
$ws.Cells.Item(111, 12).Value = 76344
$ws.Cells.Item(111, 13).Value = -7666.3334
$ws.Cells.Item(111, 14).Value = -82478

$ws.Cells.Item(122, 8).Value = 54731588
$ws.Cells.Item(122, 9).Value = 25000926
$ws.Cells.Item(122, 11).Value = 75002778
$ws.Cells.Item(122, 13).Value = -75000328

$ws.Cells.Item(132, 8).Value = 5450117.5
$ws.Cells.Item(132, 9).Value = 1322097.4
$ws.Cells.Item(132, 10).Value = 18522182
$ws.Cells.Item(132, 11).Value = 3966292.2
$ws.Cells.Item(132, 12).Value = 55566546
$ws.Cells.Item(132, 13).Value = -3963762.2
$ws.Cells.Item(132, 14).Value = -55571606

$ws.Cells.Item(138, 8).Value = 2409.74
$ws.Cells.Item(138, 9).Value = 1057.4839
$ws.Cells.Item(138, 10).Value = 3017.2754
$ws.Cells.Item(138, 11).Value = 3172.4517
$ws.Cells.Item(138, 12).Value = 9051.8262
$ws.Cells.Item(138, 13).Value = 1967.5483
$ws.Cells.Item(138, 14).Value = -19331.8262

$ws.Cells.Item(141, 8).Value = 2687.6785
$ws.Cells.Item(141, 9).Value = 734.4545000000001
$ws.Cells.Item(141, 10).Value = 9849.5
$ws.Cells.Item(141, 11).Value = 2203.3635
$ws.Cells.Item(141, 12).Value = 29548.5
$ws.Cells.Item(141, 13).Value = 2976.6365
$ws.Cells.Item(141, 14).Value = -39908.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2429493
$ws.Cells.Item(32, 9).Value = 2871884.5
$ws.Cells.Item(32, 10).Value = 47384.617
$ws.Cells.Item(32, 11).Value = 2871884.5
$ws.Cells.Item(32, 12).Value = 47384.617
$ws.Cells.Item(32, 13).Value = -2871597.5
$ws.Cells.Item(32, 14).Value = -47958.617

$ws.Cells.Item(122, 8).Value = 2446.2354
$ws.Cells.Item(122, 9).Value = 2154.182
$ws.Cells.Item(122, 10).Value = 2981.6667
$ws.Cells.Item(122, 11).Value = 6462.545999999999
$ws.Cells.Item(122, 12).Value = 8945.000100000001
$ws.Cells.Item(122, 13).Value = -4012.545999999999
$ws.Cells.Item(122, 14).Value = -13845.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1379.2593
$ws.Cells.Item(94, 9).Value = 965.2941
$ws.Cells.Item(94, 10).Value = 2083
$ws.Cells.Item(94, 11).Value = 965.2941
$ws.Cells.Item(94, 12).Value = 2083
$ws.Cells.Item(94, 13).Value = -514.2941
$ws.Cells.Item(94, 14).Value = -2985

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 13878.5
$ws.Cells.Item(50, 10).Value = 13878.5
$ws.Cells.Item(50, 12).Value = 13878.5
$ws.Cells.Item(50, 14).Value = -15128.5

$ws.Cells.Item(51, 8).Value = 33943.555
$ws.Cells.Item(51, 9).Value = 100000
$ws.Cells.Item(51, 10).Value = 15070.286
$ws.Cells.Item(51, 11).Value = 100000
$ws.Cells.Item(51, 12).Value = 15070.286
$ws.Cells.Item(51, 13).Value = -99264
$ws.Cells.Item(51, 14).Value = -16542.286

$ws.Cells.Item(60, 8).Value = 19238
$ws.Cells.Item(60, 9).Value = 28000
$ws.Cells.Item(60, 10).Value = 13761.75
$ws.Cells.Item(60, 11).Value = 28000
$ws.Cells.Item(60, 12).Value = 13761.75
$ws.Cells.Item(60, 13).Value = -27489
$ws.Cells.Item(60, 14).Value = -14783.75

$ws.Cells.Item(61, 8).Value = 33943.555
$ws.Cells.Item(61, 9).Value = 100000
$ws.Cells.Item(61, 10).Value = 15070.286
$ws.Cells.Item(61, 11).Value = 100000
$ws.Cells.Item(61, 12).Value = 15070.286
$ws.Cells.Item(61, 13).Value = -99652
$ws.Cells.Item(61, 14).Value = -15766.286

$ws.Cells.Item(68, 8).Value = 18999
$ws.Cells.Item(68, 10).Value = 18999
$ws.Cells.Item(68, 12).Value = 18999
$ws.Cells.Item(68, 14).Value = -20497

$ws.Cells.Item(71, 8).Value = 18999
$ws.Cells.Item(71, 10).Value = 18999
$ws.Cells.Item(71, 12).Value = 56997
$ws.Cells.Item(71, 14).Value = -64485

$ws.Cells.Item(74, 8).Value = 19381.334
$ws.Cells.Item(74, 10).Value = 20844.637
$ws.Cells.Item(74, 12).Value = 20844.637
$ws.Cells.Item(74, 14).Value = -22592.637

$ws.Cells.Item(77, 8).Value = 19381.334
$ws.Cells.Item(77, 10).Value = 20844.637
$ws.Cells.Item(77, 12).Value = 62533.91099999999
$ws.Cells.Item(77, 14).Value = -71269.91099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 5547.4
$ws.Cells.Item(122, 9).Value = 3662.5
$ws.Cells.Item(122, 10).Value = 7701.5713
$ws.Cells.Item(122, 11).Value = 10987.5
$ws.Cells.Item(122, 12).Value = 23104.7139
$ws.Cells.Item(122, 13).Value = -8537.5
$ws.Cells.Item(122, 14).Value = -28004.7139

$ws.Cells.Item(126, 8).Value = 6450.769
$ws.Cells.Item(126, 9).Value = 9438.923000000001
$ws.Cells.Item(126, 10).Value = 3462.6155
$ws.Cells.Item(126, 11).Value = 28316.769
$ws.Cells.Item(126, 12).Value = 10387.8465
$ws.Cells.Item(126, 13).Value = -25846.769
$ws.Cells.Item(126, 14).Value = -15327.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2554.9
$ws.Cells.Item(40, 9).Value = 2720.5715
$ws.Cells.Item(40, 10).Value = 2168.3333
$ws.Cells.Item(40, 11).Value = 2720.5715
$ws.Cells.Item(40, 12).Value = 2168.3333
$ws.Cells.Item(40, 13).Value = -2584.5715
$ws.Cells.Item(40, 14).Value = -2440.3333

$ws.Cells.Item(122, 8).Value = 8512666
$ws.Cells.Item(122, 9).Value = 1065970.2
$ws.Cells.Item(122, 10).Value = 33334984
$ws.Cells.Item(122, 11).Value = 3197910.6
$ws.Cells.Item(122, 12).Value = 100004952
$ws.Cells.Item(122, 13).Value = -3195460.6
$ws.Cells.Item(122, 14).Value = -100009852

$ws.Cells.Item(132, 8).Value = 4613375.5
$ws.Cells.Item(132, 9).Value = 5957814.5
$ws.Cells.Item(132, 10).Value = 3871.2856
$ws.Cells.Item(132, 11).Value = 17873443.5
$ws.Cells.Item(132, 12).Value = 11613.8568
$ws.Cells.Item(132, 13).Value = -17870913.5
$ws.Cells.Item(132, 14).Value = -16673.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1218.0416
$ws.Cells.Item(122, 9).Value = 1022.0625
$ws.Cells.Item(122, 10).Value = 1610
$ws.Cells.Item(122, 11).Value = 3066.1875
$ws.Cells.Item(122, 12).Value = 4830
$ws.Cells.Item(122, 13).Value = -616.1875
$ws.Cells.Item(122, 14).Value = -9730

$ws.Cells.Item(126, 8).Value = 819.625
$ws.Cells.Item(126, 9).Value = 730.3158
$ws.Cells.Item(126, 10).Value = 1159
$ws.Cells.Item(126, 11).Value = 2190.9474
$ws.Cells.Item(126, 12).Value = 3477
$ws.Cells.Item(126, 13).Value = 279.0526
$ws.Cells.Item(126, 14).Value = -8417
